$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.637.25"
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("D3").Value = "3.263.87"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "580.02"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "184.18"
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.603"
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("D9").Value = "0.129"
$ws.Range("E9").Value = "  -3.84%  "
$ws.Range("D10").Value = "6.57"
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("E11").Value = "  -3.27%  "
$ws.Range("D12").Value = "3.828.01"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").Value = "27.31"
$ws.Range("E14").Value = "  -5.48%  "
$ws.Range("D15").Value = "67.720.13"
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").Value = "0.0000167"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").Value = "3.244.91"
$ws.Range("E17").Value = "  -2.29%  "
$ws.Range("D18").Value = "5.71"
$ws.Range("E18").Value = "  -1.97%  "
$ws.Range("D19").Value = "13.42"
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("D20").Value = "398.75"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").Value = "7.55"
$ws.Range("E21").Value = "  -2.32%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "70.87"
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("D24").Value = "0.508"
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("E25").Value = "  -3.28%  "
$ws.Range("D26").Value = "0.187"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").Value = "9.52"
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("E29").Value = "  -1.96%  "
$ws.Range("D30").Value = "22.60"
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("D31").Value = "5.45"
$ws.Range("E31").Value = "  -5.08%  "
$ws.Range("D32").Value = "6.92"
$ws.Range("E32").Value = "  -3.10%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "1.24"
$ws.Range("E34").Value = "  -4.33%  "
$ws.Range("D35").Value = "163.56"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").Value = "1.45"
$ws.Range("E36").Value = "  -4.51%  "
$ws.Range("D37").Value = "1.88"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "26.88"
$ws.Range("E38").Value = "  +2.53%  "
$ws.Range("D39").Value = "0.805"
$ws.Range("E39").Value = "  -3.62%  "
$ws.Range("D40").Value = "4.50"
$ws.Range("E40").Value = "  -2.25%  "
$ws.Range("D41").Value = "2.669.39"
$ws.Range("E41").Value = "  +2.32%  "
$ws.Range("D42").Value = "6.26"
$ws.Range("E42").Value = "  -4.42%  "
$ws.Range("D43").Value = "40.71"
$ws.Range("E43").Value = "  -1.99%  "
$ws.Range("D44").Value = "0.0679"
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("D45").Value = "2.42"
$ws.Range("E45").Value = "  -5.83%  "
$ws.Range("D46").Value = "335.21"
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("D47").Value = "24.49"
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("E48").Value = "  -3.19%  "
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").Value = "0.967"
$ws.Range("E51").Value = "  -1.88%  "
